$wb = $excel.ActiveWorkbook

# --- 1. Split the existing "2022-Q3" sheet into "2022-Q4" (new data) + a
#        preserved historical "2022-Q3" copy ---------------------------------

$wsTotal = $wb.Worksheets.Item(1)          # "总计"
$wsQ3    = $wb.Worksheets.Item(2)          # "2022-Q3" (rId2 / sheetId 2)

# Rename the existing sheet in place, then duplicate it right after itself so
# the duplicate carries the OLD Q3 numbers forward under the old name, while
# the original slot becomes the new Q4 sheet we are about to fill in.
$wsQ3.Name = "2022-Q4"
$wsQ3.Copy($null, $wsQ3)

$wsOldQ3 = $wb.Worksheets.Item(3)
$wsOldQ3.Name = "2022-Q3"

$wsQ4 = $wb.Worksheets.Item(2)

# --- 2. Replace the Q4 sheet's contents with the new fund-holding data ------

$wsQ4.Cells.Clear()

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "167703"
$wsQ4.Range("C2").Value = "德邦量化优选股票（LOF）C"
$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "0.56"
$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "88.52"
$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "1.01"
$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.0057"
$wsQ4.Range("H2").Value = 5

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "167702"
$wsQ4.Range("C3").Value = "德邦量化优选股票（LOF）A"
$wsQ4.Range("D3").NumberFormat = "@"
$wsQ4.Range("D3").Value = "0.32"
$wsQ4.Range("E3").NumberFormat = "@"
$wsQ4.Range("E3").Value = "88.52"
$wsQ4.Range("F3").NumberFormat = "@"
$wsQ4.Range("F3").Value = "1.01"
$wsQ4.Range("G3").NumberFormat = "@"
$wsQ4.Range("G3").Value = "0.0032"
$wsQ4.Range("H3").Value = 5

# Data cells keep the workbook's default (unstyled) look -- clear the
# temporary "@" text format back off of them.
$wsQ4.Range("B2:G3").Style = "Normal"

# Re-apply the header / first-column formatting used elsewhere in the
# workbook (pulled from the "总计" sheet, which already carries it).
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A3").PasteSpecial(-4122)

# --- 3. Update the "总计" summary sheet -------------------------------------
# Row 2 now reports the Q4 totals; row 3 preserves the old Q3 totals.

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.01

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.54

$wsTotal.Select()
